$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the appointment/name value in A3 (was "John Test111", now "John Test1112")
$ws.Range("A3").Value = "John Test1112"

# Reflect the new active selection on the sheet (moved from C8 to A3)
$ws.Range("A3").Select()
